$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected columns to Text format so that numeric-looking
# strings (e.g. "0.996", "212.37") are stored as text, matching the
# original inline-string cell type, not auto-converted to numbers.
$cells = @(
    'D2',
    'E2',
    'D3',
    'E3',
    'D4',
    'E4',
    'D5',
    'E5',
    'D6',
    'E6',
    'D7',
    'E7',
    'D8',
    'E8',
    'D9',
    'E9',
    'D10',
    'E10',
    'D11',
    'E11',
    'D12',
    'E12',
    'D13',
    'E13',
    'D14',
    'E14',
    'D15',
    'E15',
    'D16',
    'E16',
    'D17',
    'E17',
    'D18',
    'E18',
    'B19',
    'C19',
    'D19',
    'E19',
    'B20',
    'C20',
    'D20',
    'E20',
    'D21',
    'E21',
    'D22',
    'E22',
    'D23',
    'E23',
    'D24',
    'E24',
    'D25',
    'E25',
    'D26',
    'E26',
    'B27',
    'C27',
    'D27',
    'E27',
    'B28',
    'C28',
    'D28',
    'E28',
    'D29',
    'E29',
    'D30',
    'E30',
    'D31',
    'E31',
    'D32',
    'E32',
    'D33',
    'E33',
    'D34',
    'E34',
    'E35',
    'D36',
    'E36',
    'E37',
    'D38',
    'E38',
    'D39',
    'E39',
    'B40',
    'C40',
    'D40',
    'E40',
    'D41',
    'E41',
    'B42',
    'C42',
    'D42',
    'E42',
    'B43',
    'C43',
    'D43',
    'E43',
    'D44',
    'E44',
    'D45',
    'E45',
    'D46',
    'E46',
    'D47',
    'E47',
    'D48',
    'E48',
    'B49',
    'C49',
    'D49',
    'E49',
    'B50',
    'C50',
    'D50',
    'E50',
    'D51',
    'E51',
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '27.544.29'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '1.572.52'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -1.04%  '
$ws.Range('D5').Value = '212.37'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').Value = '0.489'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('D8').Value = '22.45'
$ws.Range('E8').Value = '  +1.70%  '
$ws.Range('D9').Value = '0.251'
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').Value = '0.0596'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').Value = '0.0872'
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').Value = '1.790.03'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '1.543.83'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = '3.78'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '0.525'
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('D16').Value = '27.514.72'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').Value = '61.97'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '222.15'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '7.53'
$ws.Range('E19').Value = '  +2.12%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0710'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').Value = '0.995'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').Value = '4.18'
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').Value = '9.46'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').Value = '1.96'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').Value = '151.42'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').Value = '6.68'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '15.14'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.107'
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').Value = '1.15'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').Value = '0.0474'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('D33').Value = '1.464.30'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('D34').Value = '3.20'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('E35').Value = '  +5.66%  '
$ws.Range('D36').Value = '1.64'
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').Value = '0.544'
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.43'
$ws.Range('E40').Value = '  +4.80%  '
$ws.Range('D41').Value = '0.820'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.83'
$ws.Range('E42').Value = '  -1.20%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '0.997'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').Value = '0.973'
$ws.Range('E44').Value = '  -2.75%  '
$ws.Range('D45').Value = '65.37'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').Value = '1.79'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').Value = '1.706.63'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '86.55'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0527'
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0973'
$ws.Range('E50').Value = '  -6.41%  '
$ws.Range('D51').Value = '0.0945'
$ws.Range('E51').Value = '  -1.46%  '

# Reset style back to Normal so no stray style index is left on the cells
foreach ($addr in $cells) {
    $ws.Range($addr).Style = "Normal"
}
